# The sheet gained one new weekly price-report row. It was inserted right
# after the existing row 11 (at row 12), pushing the former rows 12-84 down
# to 13-85 (row 85 is the former row 84's data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(12).Insert()

$ws.Range("A12").Value = 2
$ws.Range("B12").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C12").Value = 'Coquimbo'
$ws.Range("D12").Value2 = 45063
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 100112026
$ws.Range("G12").Value = 'Haba'
$ws.Range("H12").Value = 'Sin especificar'
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 700
$ws.Range("K12").Value = 12000
$ws.Range("L12").Value = 14000
$ws.Range("M12").Value = 13000
$ws.Range("N12").Value = '$/malla 25 kilos'
$ws.Range("O12").Value = 'Provincia de Limarí'
$ws.Range("P12").Value = 520
$ws.Range("Q12").Value = 25
$ws.Range("R12").Value = 'Hortaliza'
